# Apply updated ticket-interest counts / price-status to the "展览" and
# "全部类型" worksheets. Both sheets contain the same rows of data and
# need identical updates.

$wb = $excel.ActiveWorkbook

# Map of row -> (F value, optional G value) changes to apply.
$updates = @(
    @{ Row = 2;  F = 364 },
    @{ Row = 4;  F = 10763 },
    @{ Row = 5;  F = 329;  G = "不可售" },
    @{ Row = 6;  F = 974 },
    @{ Row = 7;  F = 162 },
    @{ Row = 8;  F = 1332 },
    @{ Row = 9;  F = 8281 },
    @{ Row = 10; F = 37 },
    @{ Row = 12; F = 341 },
    @{ Row = 15; F = 3295 },
    @{ Row = 17; F = 326 },
    @{ Row = 18; F = 20 },
    @{ Row = 19; F = 778 },
    @{ Row = 21; F = 1068 },
    @{ Row = 23; F = 106 },
    @{ Row = 24; F = 1764 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Range("G" + $u.Row).Value = $u.G
        }
    }
}
